$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60. This pushes the existing rows 60:90
# down to 61:91 (preserving all of their data/styles), and the new,
# empty row 60 inherits formatting (e.g. the date style on column D)
# from the row above it.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new weekly record. The
# "constant" columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat the same market /
# product metadata used throughout this block of rows.
$ws.Range("A60").Value = 5
$ws.Range("B60").Value = "Macroferia Regional de Talca"
$ws.Range("C60").Value = "Maule"
$ws.Range("D60").Value = 44572
$ws.Range("E60").Value = 7
$ws.Range("F60").Value = 100112030
$ws.Range("G60").Value = "Poroto granado"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 23000
$ws.Range("L60").Value = 23000
$ws.Range("M60").Value = 23000
$ws.Range("N60").Value = "$/saco 25 kilos"
$ws.Range("O60").Value = "Región del Maule"
$ws.Range("P60").Value = 920
$ws.Range("Q60").Value = 25
$ws.Range("R60").Value = "Hortaliza"
